$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.651.63'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.818.42'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.582'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.81%  '
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '34.64'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.301'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0701'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0953'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.082.18'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.39'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.817.72'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.644'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.664.58'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.33'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.17'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0802'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '246.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.55'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.26%  '
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.18'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '173.73'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.47%  '
$ws.Range("E25").Value = '  +1.20%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.49'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.80'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.27%  '
$ws.Range("E28").Value = '  +2.71%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.99'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0531'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.84'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.44%  '
$ws.Range("E33").Value = '  +1.10%  '
$ws.Range("E34").Value = '  +0.75%  '
$ws.Range("E35").Value = '  +1.41%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.406.77'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.67%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.680'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.54%  '
$ws.Range("E38").Value = '  +1.73%  '
$ws.Range("E39").Value = '  +0.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '84.37'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.88'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.950'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.73%  '
$ws.Range("E43").Value = '  -0.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.72'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.67%  '
$ws.Range("E45").Value = '  +2.85%  '
$ws.Range("E46").Value = '  -1.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.07'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.983.43'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.48'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.24%  '
$ws.Range("E50").Value = '  -1.35%  '
$ws.Range("E51").Value = '  +0.18%  '
